$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 4.5
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 2
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 2.37
$ws.Range("AP2").Value = 9.5

# Row 3 updates
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 4.33
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.48
$ws.Range("AA3").Value = 2.1
$ws.Range("AB3").Value = 1.67
$ws.Range("AJ3").Value = 6
$ws.Range("AK3").Value = 19
$ws.Range("AL3").Value = 67
$ws.Range("AN3").Value = 9.5
